$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('126 Racecourse Road Public Housing Tower Flemington', 6),
    @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 18),
    @('Al Haj Halal Meats Glenroy', 11),
    @('Al-Taqwa College Truganina', 8),
    @('Apartment Complex 120 Racecourse Road North Melbourne', 5),
    @('Australia Post Distribution Centre Sunshine West', 5),
    @('Australian Lamb Colac East', 5),
    @('Baker Bleu Caulfield North', 7),
    @('Baxter Foods Australia Campbellfield', 5),
    @('CFMEU Melbourne Office', 5),
    @('CS Square Caroline Springs', 13),
    @('Cafe Roco Dandenong', 6),
    @('Campbellfield Ford Complex Vaccination Clinic Campbellfield', 10),
    @('Cardinia Lakes Early Learning Centre Pakenham', 5),
    @('Carlton Finishing Pty. Ltd. Campbellfield', 12),
    @('Caroline Springs Police Station', 8),
    @('Chemist Warehouse Campbellfield DC', 5),
    @('Chemist Warehouse Fillo Drive Somerton', 9),
    @('City of Wyndham Community', 6),
    @('Coles Campbellfield Plaza Campbellfield', 7),
    @('Coles Coburg North Village', 18),
    @('Coles Pakenham Place Shopping Centre', 6),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 16),
    @('Construction Site 1 Warde Street Footscray', 5),
    @('Construction Site Olea Apartment Caulfield North', 11),
    @('Costco Wholesale Epping', 26),
    @('Crusader Caravans Epping', 22),
    @('Dandenong Police Station Dandenong', 5),
    @('DayHab Rehabilitation Treatment Centre Ringwood East', 7),
    @('Direct Freight Express Campbellfield', 5),
    @('Disability Residence Life without Barriers Ashwood', 5),
    @('Don Watson Coldstore Derrimutg', 5),
    @('Epworth Healthcare Epworth Richmond Emergency Department', 6),
    @('Ermha365 Residential Disability Care Services Paperbark St Doveton', 9),
    @('FedEx Station Melbourne Airport', 15),
    @('Fine Food Holdings Pty Ltd Dandenong South', 10),
    @('Fitzroy St Precinct St Kilda', 7),
    @('Fonterra Manufacturing Workplace Campbellfield', 9),
    @('General Foods Campbellfield', 14),
    @('Gladstone Parade Early Learning & Kinder Glenroy', 7),
    @('Goodstart Early Learning Altona', 9),
    @('Green Leaves Early Leaning Cairnlea', 5),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 15),
    @('Greenvale Primary School', 5),
    @('HEI Schools Emerald Early Learning Centre Emerald', 5),
    @('Hamilton Marino 236 Jasper Road McKinnon', 9),
    @('Hello Fresh Warehouse Ravenhall', 7),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 9),
    @('Industrial Galvanizers Valmont Coatings Campbellfield', 14),
    @('Inghams Enterprises Thomastown', 6),
    @('Kippers Seafood Werribee', 6),
    @('Kool Kidz Childcare Narre Warren', 13),
    @('Lineage Logistics Laverton North', 6),
    @('Linfox Somerton National Distribution Centre Somerton', 9),
    @('Mecca Distribution Centre Warehouse Melbourne Airport', 9),
    @('Melbourne Assessment Prison West Melbourne', 9),
    @('Melbourne Metropolitan Remand Centre Ravenhall', 7),
    @('Melbourne West Police Station Docklands', 6),
    @('Mill Park Police Station Mill Park', 8),
    @('MyCentre Childcare Broadmeadows', 21),
    @('National Gallery of Victoria Melbourne', 7),
    @('Nido Early School Ascot Vale', 12),
    @('Nido Early School Glenroy', 23),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 49),
    @('Northern Health The Northern Hospital Epping', 16),
    @('OnQ Plumbing and Excavations Craigieburn', 13),
    @('Oporto Coolaroo', 6),
    @('Oscar Romero Catholic Primary School Craigieburn', 5),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 9),
    @('Pacific Meat Thomastown', 5),
    @('Private Residence Daycare Allumba Way Wollert', 8),
    @('Ravenhall Correctional Centre Ravenhall', 10),
    @('Richmond Quarter 261-271 Bridge Road Construction Site Richmond', 11),
    @('Sacca''s Fruit World Broadmeadows Central Shopping Centre', 6),
    @('Sharpline Stainless Steel Coburg North', 6),
    @('St Margaret''s Primary School OSHC Maribyrnong', 12),
    @('St Vincents Hospital Emergency Department Melbourne', 11),
    @('Tek Foods Somerton', 5),
    @('The Huntly-Goornong Rail Works', 6),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 16),
    @('The Royal Melbourne Hospital Parkville', 6),
    @('The Royal Melbourne Hospital Parkville Emergency Department', 6),
    @('The Royal Melbourne Hospital Ward 6SE Parkville', 20),
    @('The Royal Talbot Rehabilitation Centre Kew', 11),
    @('ThorwestenCabinets Pakenham', 11),
    @('Truganina Early Learning Centre Truganina', 7),
    @('Visy Recycling Springvale', 22),
    @('Wallaby Childcare Wollert', 18),
    @('Werribee Mercy Hospital Emergency Department', 15),
    @('Western Health Footscray Hospital Emergency Department', 9),
    @('Western Health Sunshine Hospital Emergency Department', 12),
    @('Western Health Sunshine Hospital GEM Ward St Albans', 6),
    @('Yarra Childcare Centre Truganina', 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

